# Fruta / hortaliza, semanal
# Rotate the data (excluding A,B,C,E,F,G,H,I,J,K,Q,T which are unchanged)
# for rows 4,5,6: new row4 = old row6, new row5 = old row4, new row6 = old row5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that change, for rows 4, 5 and 6
$cols = @("D","L","M","N","O","P","R","S")

$orig4 = @{}
$orig5 = @{}
$orig6 = @{}
foreach ($col in $cols) {
    $orig4[$col] = $ws.Range("${col}4").Value2
    $orig5[$col] = $ws.Range("${col}5").Value2
    $orig6[$col] = $ws.Range("${col}6").Value2
}

# Apply rotation: row4 <- old row6, row5 <- old row4, row6 <- old row5
foreach ($col in $cols) {
    $ws.Range("${col}4").Value = $orig6[$col]
    $ws.Range("${col}5").Value = $orig4[$col]
    $ws.Range("${col}6").Value = $orig5[$col]
}
